# `excel`: minor performance tweaks — also added test for "data types" sheet.
$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet (date_test) and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "data types"

# Fill column F (string) top-to-bottom first, then column G (emojis),
# then the header row, then column H (foreign) -- this reproduces the
# exact shared-string insertion order of the authored workbook.
$ws.Range("F2").Value = "The"
$ws.Range("F3").Value = "quick"
$ws.Range("F4").Value = "brown"
$ws.Range("F5").Value = "fox"
$ws.Range("F6").Value = "jumped"

$ws.Range("G2").Value = "The"
$ws.Range("G3").Value = "🍔"
$ws.Range("G4").Value = "is"
$ws.Range("G5").Value = "💩"
$ws.Range("G6").Value = "🙀"

$ws.Range("A1").Value = "int"
$ws.Range("B1").Value = "float"
$ws.Range("C1").Value = "bool"
$ws.Range("D1").Value = "date"
$ws.Range("E1").Value = "duration"
$ws.Range("F1").Value = "string"
$ws.Range("G1").Value = "emojis"
$ws.Range("H1").Value = "foreign"

$ws.Range("H2").Value = "敏捷的棕色狐狸在森林里奔跑"
$ws.Range("H3").Value = "Franz jagt im komplett verwahrlosten Taxi quer durch Bayern"
$ws.Range("H4").Value = "Le rusé goupil franchit d'un bond le chien somnolent."
$ws.Range("H5").Value = "El rápido zorro marrón"
$ws.Range("H6").Value = "いろはにほへとちりぬるをわかよたれそつねならむうゐのおくやまけふこえてあさきゆめみしゑひもせす"

# Column A -- int
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Column B -- float
$ws.Range("B2").Value = 1.1
$ws.Range("B3").Value = 1.32434354545454
$ws.Range("B4").Value = 0.42354645656453399
$ws.Range("B5").Value = -54545.656575678498
$ws.Range("B6").Value = -5446563454.4354601

# Column C -- bool (rows 2,3,6 are real booleans; rows 4,5 are plain numbers)
$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $false
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = $true

# Column D -- date (rows 2-4 date only, rows 5-6 date+time)
$ws.Range("D2").Value = 37145
$ws.Range("D2").NumberFormat = "m/d/yy"
$ws.Range("D3").Value = 45206
$ws.Range("D3").NumberFormat = "m/d/yy"
$ws.Range("D4").Value = 15317
$ws.Range("D4").NumberFormat = "m/d/yy"
$ws.Range("D5").Value = 37145.354166666664
$ws.Range("D5").NumberFormat = "m/d/yy h:mm"
$ws.Range("D6").Value = 16655.34375
$ws.Range("D6").NumberFormat = "m/d/yy h:mm"

# Column E -- duration ([h]:mm:ss number format -> numFmtId 46)
$ws.Range("E2").Value = 0.43055555555555558
$ws.Range("E2").NumberFormat = "[h]:mm:ss"
$ws.Range("E3").Value = 0.98984953703703704
$ws.Range("E3").NumberFormat = "[h]:mm:ss"
$ws.Range("E4").Value = 1.2815162037037038
$ws.Range("E4").NumberFormat = "[h]:mm:ss"
$ws.Range("E5").Value = 0.97916666666666663
$ws.Range("E5").NumberFormat = "[h]:mm:ss"
$ws.Range("E6").Value = 0.00046296296296296293
$ws.Range("E6").NumberFormat = "[h]:mm:ss"

# Match the authored column best-fit widths for columns B and D.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).AutoFit()

$ws.Range("E3").Select()
